$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string into a cell without letting Excel
# auto-coerce numeric-looking text (e.g. "607.86") into a Number.
# Builds the text via a `="..."` formula in a scratch cell, copies it,
# and pastes-special "values only" into the destination - this bakes in
# a plain text value (like a formula->value paste in real Excel) without
# requiring a quote-prefix or a Text number-format (so no style changes).
$scratch = $ws.Range("Z1")
function Set-TextValue($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") '69.868.37'

Set-TextValue $ws.Range("D3") '3.506.00'
$ws.Range("E3").Value = '  -1.60%  '

$ws.Range("E4").Value = '  -0.12%  '

Set-TextValue $ws.Range("D5") '607.86'
$ws.Range("E5").Value = '  +3.42%  '

Set-TextValue $ws.Range("D6") '191.91'
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("E7").Value = '  +0.74%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("E9").Value = '  -0.13%  '

Set-TextValue $ws.Range("D10") '0.664'
$ws.Range("E10").Value = '  +2.93%  '

$ws.Range("E11").Value = '  -1.19%  '

$ws.Range("E12").Value = '  -0.76%  '

$ws.Range("E13").Value = '  +2.22%  '

Set-TextValue $ws.Range("D14") '4.064.30'
$ws.Range("E14").Value = '  -1.49%  '

Set-TextValue $ws.Range("D15") '619.57'
$ws.Range("E15").Value = '  +10.42%  '

Set-TextValue $ws.Range("D16") '69.941.71'
$ws.Range("E16").Value = '  -1.13%  '

Set-TextValue $ws.Range("D18") '18.90'
$ws.Range("E18").Value = '  -0.25%  '

Set-TextValue $ws.Range("D19") '3.515.03'
$ws.Range("E19").Value = '  -2.86%  '

$ws.Range("E20").Value = '  -0.15%  '

$ws.Range("E21").Value = '  -0.37%  '

Set-TextValue $ws.Range("D22") '17.65'
$ws.Range("E22").Value = '  -1.47%  '

Set-TextValue $ws.Range("D23") '105.94'
$ws.Range("E23").Value = '  +12.96%  '

$ws.Range("E24").Value = '  +0.66%  '

Set-TextValue $ws.Range("D25") '5.00'
$ws.Range("E25").Value = '  +2.04%  '

$ws.Range("E26").Value = '  +3.92%  '

$ws.Range("E27").Value = '  -1.01%  '

Set-TextValue $ws.Range("D28") '9.85'
$ws.Range("E28").Value = '  +5.71%  '

Set-TextValue $ws.Range("D29") '34.14'
$ws.Range("E29").Value = '  +5.58%  '

$ws.Range("E30").Value = '  +0.75%  '

$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D31") '12.64'
$ws.Range("E31").Value = '  +3.69%  '

$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D32") '4.15'
$ws.Range("E32").Value = '  +4.71%  '

$ws.Range("E33").Value = '  -0.11%  '

Set-TextValue $ws.Range("D34") '64.22'
$ws.Range("E34").Value = '  +1.68%  '

Set-TextValue $ws.Range("D35") '3.721.04'
$ws.Range("E35").Value = '  +1.88%  '

$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D36") '3.10'
$ws.Range("E36").Value = '  -4.23%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D37") '525.12'
$ws.Range("E37").Value = '  -0.84%  '

$ws.Range("E38").Value = '  -0.05%  '

Set-TextValue $ws.Range("D39") '0.0₃0798'
$ws.Range("E39").Value = '  +1.54%  '

Set-TextValue $ws.Range("D40") '0.392'
$ws.Range("E40").Value = '  -3.59%  '

$ws.Range("E41").Value = '  -3.49%  '

Set-TextValue $ws.Range("D42") '3.57'
$ws.Range("E42").Value = '  +0.28%  '

$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("E44").Value = '  +1.07%  '

Set-TextValue $ws.Range("D45") '2.87'
$ws.Range("E45").Value = '  -2.50%  '

$ws.Range("E46").Value = '  +2.58%  '

Set-TextValue $ws.Range("D47") '3.33'
$ws.Range("E47").Value = '  -4.00%  '

$ws.Range("E48").Value = '  -4.66%  '

$ws.Range("E49").Value = '  +0.48%  '

Set-TextValue $ws.Range("D50") '133.14'
$ws.Range("E50").Value = '  -1.74%  '

$ws.Range("E51").Value = '  -7.17%  '

$scratch.ClearContents()